$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.19"
$ws.Range("D3").Value = "'23.07"
$ws.Range("D4").Value = "'5.405"
$ws.Range("D5").Value = "'0.05923"
$ws.Range("D6").Value = "'3.450"
$ws.Range("D7").Value = "'6.541"
$ws.Range("D9").Value = "'0.9120"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("D11").Value = "'0.07350"
$ws.Range("D12").Value = "'0.03265"
$ws.Range("D13").Value = "'0.03042"
$ws.Range("D14").Value = "'0.09352"
$ws.Range("D15").Value = "'3.852"
$ws.Range("D16").Value = "'0.001579"
$ws.Range("D17").Value = "'0.04680"
$ws.Range("D18").Value = "'0.01119"
$ws.Range("E18").Value = "17OneONEBestin24h"
$ws.Range("D19").Value = "'0.006124"
$ws.Range("D20").Value = "'0.004976"
$ws.Range("D21").Value = "'0.0009808"
$ws.Range("D22").Value = "'0.00009405"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("D23").Value = "'3.610"
$ws.Range("D27").Value = "'0.0002902"
$ws.Range("D41").Value = "'0.006198"
$ws.Range("D42").Value = "'0.1076"
$ws.Range("D43").Value = "'0.003002"
$ws.Range("D44").Value = "'0.008203"
$ws.Range("D45").Value = "'0.00005246"
$ws.Range("D47").Value = "'0.7824"
$ws.Range("D48").Value = "'0.002271"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
